$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format Price column as text first so numeric-looking values are not
# auto-converted by Excel, then reset the style so no extra style index
# is left behind on the cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '36.427.25'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '1.939.56'
$ws.Range("E3").Value = '  -2.08%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '242.48'
$ws.Range("E5").Value = '  -1.39%  '
$ws.Range("D6").Value = '0.609'
$ws.Range("E6").Value = '  -2.47%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '57.00'
$ws.Range("E8").Value = '  -3.78%  '
$ws.Range("D9").Value = '0.358'
$ws.Range("E9").Value = '  -4.37%  '
$ws.Range("E10").Value = '  -2.93%  '
$ws.Range("D11").Value = '0.102'
$ws.Range("E11").Value = '  -1.94%  '
$ws.Range("D12").Value = '2.224.66'
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("D13").Value = '21.25'
$ws.Range("E13").Value = '  -5.37%  '
$ws.Range("E14").Value = '  -5.89%  '
$ws.Range("E15").Value = '  -3.05%  '
$ws.Range("D16").Value = '5.13'
$ws.Range("E16").Value = '  -6.24%  '
$ws.Range("D17").Value = '1.930.63'
$ws.Range("E17").Value = '  -2.20%  '
$ws.Range("D18").Value = '36.386.52'
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("D19").Value = '69.10'
$ws.Range("E19").Value = '  -1.97%  '
$ws.Range("D20").Value = '0.0₃0863'
$ws.Range("E20").Value = '  -4.80%  '
$ws.Range("D21").Value = '227.50'
$ws.Range("E21").Value = '  -3.08%  '
$ws.Range("D22").Value = '4.97'
$ws.Range("E22").Value = '  -5.85%  '
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("E24").Value = '  -6.70%  '
$ws.Range("D25").Value = '2.28'
$ws.Range("E25").Value = '  -1.28%  '
$ws.Range("D26").Value = '9.19'
$ws.Range("E26").Value = '  -6.61%  '
$ws.Range("D27").Value = '161.19'
$ws.Range("E27").Value = '  -2.54%  '
$ws.Range("D28").Value = '0.132'
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("D29").Value = '19.20'
$ws.Range("E29").Value = '  -3.75%  '
$ws.Range("D30").Value = '0.118'
$ws.Range("E30").Value = '  -1.95%  '
$ws.Range("D31").Value = '1.09'
$ws.Range("E31").Value = '  -7.82%  '
$ws.Range("D32").Value = '4.55'
$ws.Range("E32").Value = '  -7.01%  '
$ws.Range("E33").Value = '  -5.23%  '
$ws.Range("D34").Value = '4.16'
$ws.Range("E34").Value = '  -6.19%  '
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("D36").Value = '6.07'
$ws.Range("E36").Value = '  -0.36%  '
$ws.Range("E37").Value = '  -0.80%  '
$ws.Range("D38").Value = '2.17'
$ws.Range("E38").Value = '  -2.47%  '
$ws.Range("E39").Value = '  +5.87%  '
$ws.Range("D40").Value = '0.0988'
$ws.Range("E40").Value = '  +2.29%  '
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("E42").Value = '  -2.18%  '
$ws.Range("E43").Value = '  -5.37%  '
$ws.Range("D44").Value = '15.61'
$ws.Range("E44").Value = '  -4.17%  '
$ws.Range("D45").Value = '1.339.88'
$ws.Range("E45").Value = '  -1.96%  '
$ws.Range("D47").Value = '86.30'
$ws.Range("E47").Value = '  -5.50%  '
$ws.Range("E48").Value = '  -4.38%  '
$ws.Range("D49").Value = '2.83'
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("D50").Value = '2.115.92'
$ws.Range("E50").Value = '  -1.95%  '
$ws.Range("D51").Value = '43.01'
$ws.Range("E51").Value = '  -5.38%  '

$ws.Range("D2:D51").Style = "Normal"
